# Edit AggTagTemplate.xlsx: switch the California/Nevada "groupBy" attribute
# from "nameFirstLetter" to "name.substring(0,1)", and set the active
# selection on the other sheets to A2 (the first data row).

$wb = $excel.ActiveWorkbook

# --- 1. Update the two jt:agg cells on the "Agg" sheet -----------------
$aggSheet = $wb.Worksheets.Item("Agg")

$caText = '<jt:agg items="${california.counties}" aggs="Count(*);Sum(population);StdDev(area)" aggsVar="aggs" valuesVar="values" groupBy="name.substring(0,1)"><jt:forEach items="${values}" var="value">${value.object.nameFirstLetter}'
$nvText = '<jt:agg items="${nevada.counties}" aggs="Count(*);Sum(population);StdDev(area)" aggsVar="aggs" valuesVar="values" groupBy="name.substring(0,1)"><jt:forEach items="${values}" var="value">${value.object.nameFirstLetter}'

$aggSheet.Range("A3").Value = $caText
$aggSheet.Range("A6").Value = $nvText

# --- 2. Move the selection to A2 on the other data sheets --------------
$wb.Worksheets.Item("Msd").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Rollup").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Rollups").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("GroupingSets").Range("A2").Select() | Out-Null

# Restore the originally active sheet/selection ("Agg", A1:D1).
$aggSheet.Activate()
